# Auto-generated edit script: apply numeric corrections to Leve profit tables
# across multiple worksheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR), per the commit diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2510.2354
$ws.Range("J62").Value = 2615.8
$ws.Range("L62").Value = 2615.8
$ws.Range("N62").Value = -3863.8
# Row 64
$ws.Range("H64").Value = 3002.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3002.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3002.5
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3498.5
# Row 65
$ws.Range("H65").Value = 2510.2354
$ws.Range("J65").Value = 2615.8
$ws.Range("L65").Value = 13079
$ws.Range("N65").Value = -19319
# Row 67
$ws.Range("H67").Value = 3002.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3002.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3002.5
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4718.5
# Row 112
$ws.Range("H112").Value = 1520.0625
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1588.0667
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 4764.2001
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -6980.2001
# Row 129
$ws.Range("H129").Value = 429515
$ws.Range("I129").Value = 538.7143
$ws.Range("J129").Value = 471808.44
$ws.Range("K129").Value = 1616.1429
$ws.Range("L129").Value = 1415425.32
$ws.Range("M129").Value = 3383.8571
$ws.Range("N129").Value = -1425425.32
# Row 138
$ws.Range("H138").Value = 1700.0652
$ws.Range("I138").Value = 1040.1072
$ws.Range("J138").Value = 2726.6667
$ws.Range("K138").Value = 3120.3216
$ws.Range("L138").Value = 8180.000100000001
$ws.Range("M138").Value = 2019.6784
$ws.Range("N138").Value = -18460.0001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17881.045
$ws.Range("I32").Value = 20675.176
$ws.Range("J32").Value = 3402.3635
$ws.Range("K32").Value = 20675.176
$ws.Range("L32").Value = 3402.3635
$ws.Range("M32").Value = -20388.176
$ws.Range("N32").Value = -3976.3635
# Row 45
$ws.Range("H45").Value = 960.5
$ws.Range("I45").Value = 974.5
$ws.Range("J45").Value = 953.5
$ws.Range("K45").Value = 974.5
$ws.Range("L45").Value = 953.5
$ws.Range("M45").Value = -597.5
$ws.Range("N45").Value = -1707.5
# Row 74
$ws.Range("H74").Value = 1238
$ws.Range("I74").Value = 1800
$ws.Range("J74").Value = 957
$ws.Range("K74").Value = 1800
$ws.Range("L74").Value = 957
$ws.Range("M74").Value = -926
$ws.Range("N74").Value = -2705
# Row 77
$ws.Range("H77").Value = 1238
$ws.Range("I77").Value = 1800
$ws.Range("J77").Value = 957
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 4785
$ws.Range("M77").Value = -4632
$ws.Range("N77").Value = -13521
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 132
$ws.Range("H132").Value = 5082.6587
$ws.Range("I132").Value = 7527
$ws.Range("J132").Value = 3351.25
$ws.Range("K132").Value = 22581
$ws.Range("L132").Value = 10053.75
$ws.Range("M132").Value = -20051
$ws.Range("N132").Value = -15113.75

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 141
$ws.Range("H141").Value = 42445
$ws.Range("J141").Value = 42445
$ws.Range("L141").Value = 42445
$ws.Range("N141").Value = -52805

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 349.83334
$ws.Range("I22").Value = 336.18182
$ws.Range("K22").Value = 336.18182
$ws.Range("M22").Value = 13.81817999999998
# Row 31
$ws.Range("H31").Value = 3775527.5
$ws.Range("I31").Value = 3384.2
$ws.Range("J31").Value = 5264531.5
$ws.Range("K31").Value = 3384.2
$ws.Range("L31").Value = 5264531.5
$ws.Range("M31").Value = -3089.2
$ws.Range("N31").Value = -5265121.5
# Row 34
$ws.Range("H34").Value = 3775527.5
$ws.Range("I34").Value = 3384.2
$ws.Range("J34").Value = 5264531.5
$ws.Range("K34").Value = 3384.2
$ws.Range("L34").Value = 5264531.5
$ws.Range("M34").Value = -3182.2
$ws.Range("N34").Value = -5264935.5
# Row 58
$ws.Range("H58").Value = 9804884
$ws.Range("I58").Value = 878.7308
$ws.Range("J58").Value = 41667900
$ws.Range("K58").Value = 878.7308
$ws.Range("L58").Value = 41667900
$ws.Range("M58").Value = -675.7308
$ws.Range("N58").Value = -41668306
# Row 132
$ws.Range("H132").Value = 2584.9412
$ws.Range("I132").Value = 2308.8
$ws.Range("K132").Value = 6926.400000000001
$ws.Range("M132").Value = -4396.400000000001
# Row 136
$ws.Range("H136").Value = 9804884
$ws.Range("I136").Value = 878.7308
$ws.Range("J136").Value = 41667900
$ws.Range("K136").Value = 2636.1924
$ws.Range("L136").Value = 125003700
$ws.Range("M136").Value = -86.19239999999991
$ws.Range("N136").Value = -125008800

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 18000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# Row 58
$ws.Range("H58").Value = 250006130
$ws.Range("I58").Value = 1000000000
$ws.Range("J58").Value = 8166.6665
$ws.Range("K58").Value = 1000000000
$ws.Range("L58").Value = 8166.6665
$ws.Range("M58").Value = -999999723
$ws.Range("N58").Value = -8720.666499999999
# Row 70
$ws.Range("H70").Value = 9111777
$ws.Range("I70").Value = 11864687
$ws.Range("K70").Value = 11864687
$ws.Range("M70").Value = -11864417
# Row 73
$ws.Range("H73").Value = 9111777
$ws.Range("I73").Value = 11864687
$ws.Range("K73").Value = 11864687
$ws.Range("M73").Value = -11863751
# Row 122
$ws.Range("H122").Value = 50004820
$ws.Range("I122").Value = 100006850
$ws.Range("J122").Value = 2790
$ws.Range("K122").Value = 300020550
$ws.Range("L122").Value = 8370
$ws.Range("M122").Value = -300018100
$ws.Range("N122").Value = -13270
# Row 126
$ws.Range("H126").Value = 2617.5293
$ws.Range("I126").Value = 914
$ws.Range("J126").Value = 3810
$ws.Range("K126").Value = 2742
$ws.Range("L126").Value = 11430
$ws.Range("M126").Value = -272
$ws.Range("N126").Value = -16370
# Row 132
$ws.Range("H132").Value = 77089.19
$ws.Range("I132").Value = 102434.2
$ws.Range("J132").Value = 4674.857
$ws.Range("K132").Value = 307302.6
$ws.Range("L132").Value = 14024.571
$ws.Range("M132").Value = -304772.6
$ws.Range("N132").Value = -19084.571

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1873.0476
$ws.Range("I40").Value = 1711.25
$ws.Range("J40").Value = 2390.8
$ws.Range("K40").Value = 1711.25
$ws.Range("L40").Value = 2390.8
$ws.Range("M40").Value = -1575.25
$ws.Range("N40").Value = -2662.8
# Row 55
$ws.Range("H55").Value = 280.89474
$ws.Range("I55").Value = 253.35715
$ws.Range("J55").Value = 358
$ws.Range("K55").Value = 253.35715
$ws.Range("L55").Value = 358
$ws.Range("M55").Value = -80.35714999999999
$ws.Range("N55").Value = -704
# Row 119
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
# Row 122
$ws.Range("H122").Value = 3835.5
$ws.Range("I122").Value = 7267.3335
$ws.Range("J122").Value = 2899.5454
$ws.Range("K122").Value = 21802.0005
$ws.Range("L122").Value = 8698.636200000001
$ws.Range("M122").Value = -19352.0005
$ws.Range("N122").Value = -13598.6362
# Row 132
$ws.Range("H132").Value = 10233.037
$ws.Range("I132").Value = 12573.723
$ws.Range("J132").Value = 5551.6665
$ws.Range("K132").Value = 37721.169
$ws.Range("L132").Value = 16654.9995
$ws.Range("M132").Value = -35191.169
$ws.Range("N132").Value = -21714.9995

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2791.6667
$ws.Range("I122").Value = 2950
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8850
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6400
$ws.Range("N122").Value = -10900
# Row 136
$ws.Range("H136").Value = 2152.04
$ws.Range("I136").Value = 2585.0715
$ws.Range("J136").Value = 1600.909
$ws.Range("K136").Value = 7755.2145
$ws.Range("L136").Value = 4802.727000000001
$ws.Range("M136").Value = -5205.2145
$ws.Range("N136").Value = -9902.727000000001
